# Apply the "نواقص الأصناف" (shortage items) report update:
#  - insert a new product row ("بادي سبلاش ايفا" / 175.00) in its alphabetically
#    sorted place (pushing every row below it down by one)
#  - drop the "لزق بثور" row (it shifts out at the bottom of that block)
#  - refresh the grand-total cell and the generated timestamp footer
#
# The sheet's row/merge/style layout does not change -- only the text & number
# values carried by each row shift. So instead of doing a literal
# insert/delete of worksheet rows (which would drag merged-cell / style
# metadata along in ways that are hard to control from COM), we simply
# rewrite the value of every affected cell to match the new, shifted table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column map for the product table (rows 7-30):
#   C = item name (merged C:G)
#   H = current balance ratio (merged H:K)
#   N = price (merged N:O)
#   P = selling price  (numeric display format, but the source keeps it as
#       literal text so the 4-decimal value survives the "0.00" format --
#       we have to round-trip NumberFormat to force text storage too)
#   Q = number of deals ratio (single cell)

$rows = @(
    @{ Row = 22; C = "بادي سبلاش ايفا";        H = "1:0"; N = "175.00"; P = "175.0000"; Q = "1:0" },
    @{ Row = 23; C = "بادي لوشن كير اند مور";   H = "2:0"; N = "85.00";  P = "85.0000";  Q = "1:0" },
    @{ Row = 24; C = "بيرسول حشرات طائره";      H = "0:0"; N = "50.00";  P = "50.0000";  Q = "1:0" },
    @{ Row = 25; C = "حبايه";                   H = "0:0"; N = "3.00";   P = "3.0000";   Q = "1:0" },
    @{ Row = 26; C = "زيت فاتيكا وسط 90 مل";    H = "9:0"; N = "25.00";  P = "25.0000";  Q = "1:0" },
    @{ Row = 27; C = "سرنجات 3 سم";             H = "0:0"; N = "2.00";   P = "14.0000";  Q = "7:0" },
    @{ Row = 28; C = "شمع حريمي";               H = "8:0"; N = "50.00";  P = "50.0000";  Q = "1:0" },
    @{ Row = 29; C = "كالونا ";                 H = "0:0"; N = "15.00";  P = "15.0000";  Q = "1:0" },
    @{ Row = 30; C = "معطر جو FRIDA ";          H = "8:0"; N = "65.00";  P = "65.0000";  Q = "1:0" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C    # C - item name
    $ws.Cells.Item($r.Row, 8).Value = $r.H    # H - current balance ratio
    $ws.Cells.Item($r.Row, 14).Value = $r.N   # N - price

    # P - selling price: cell is numeric-formatted ("0.00") but must keep the
    # full literal text (e.g. "175.0000") instead of being parsed into a
    # rounded number. Flip to a text format while writing, then restore the
    # original numeric format so the style index is unchanged.
    $pCell = $ws.Cells.Item($r.Row, 16)
    $pCell.NumberFormat = "@"
    $pCell.Value = $r.P
    $pCell.NumberFormat = "0.00"

    $ws.Cells.Item($r.Row, 17).Value = $r.Q   # Q - deals ratio
}

# Grand total (sum of the P column) picks up the new item's price and loses
# the removed item's price: 1184.495 + 175.00 - 85.00 = 1274.495
$ws.Cells.Item(31, 16).Value = 1274.4949999999999

# Footer timestamp: regenerated at a later time of day
$ws.Cells.Item(32, 1).Value = "Friday, 19 September, 2025 4:11 PM"
